$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 214, shifting existing rows 214-217 down to 215-218.
$ws.Rows.Item(214).Insert()

# Populate the newly inserted row 214 with the new weekly record.
$ws.Range("A214").Value = 7
$ws.Range("B214").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C214").Value = "Ñuble"
$ws.Range("D214").Value = 44595
$ws.Range("D214").NumberFormat = $ws.Range("D215").NumberFormat
$ws.Range("E214").Value = 16
$ws.Range("F214").Value = 100112002
$ws.Range("G214").Value = "Pimiento"
$ws.Range("H214").Value = "Cuatro cascos verde"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 120
$ws.Range("K214").Value = 6500
$ws.Range("L214").Value = 7000
$ws.Range("M214").Value = 6750
$ws.Range("N214").Value = "`$/caja 15 kilos"
$ws.Range("O214").Value = "Región del Maule"
$ws.Range("P214").Value = 450
$ws.Range("Q214").Value = 15
$ws.Range("R214").Value = "Hortaliza"
